$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds date-like strings (e.g. "20201211") that must stay text,
# not get auto-converted to numbers. Format as Text first so the values
# round-trip as inline strings.
$ws.Range("A118:A121").NumberFormat = "@"

# Update existing row 118
$ws.Range("A118").Value = "20201211"
$ws.Range("B118").Value = 605.0

# Add new row 119
$ws.Range("A119").Value = "20201222"
$ws.Range("B119").Value = 653.0

# Add new row 120
$ws.Range("A120").Value = "20210108"
$ws.Range("B120").Value = 539.0

# Add new row 121
$ws.Range("A121").Value = "20210329"
$ws.Range("B121").Value = 736.0
